# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows: 3, 4, 5, 8, 10)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 96
$ws1.Range("F4").Value = 490
$ws1.Range("F5").Value = 4866
$ws1.Range("F8").Value = 291
$ws1.Range("F10").Value = 226

# Sheet "全部类型" (rows: 3, 4, 5, 8, 11)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 96
$ws4.Range("F4").Value = 490
$ws4.Range("F5").Value = 4866
$ws4.Range("F8").Value = 291
$ws4.Range("F11").Value = 226
